$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33:G33").Copy() | Out-Null
$ws.Range("A34:G35").Insert(-4121) | Out-Null

$ws.Range("D39:D40").Cut($ws.Range("D37")) | Out-Null
$ws.Range("A39:G40").Delete(-4162) | Out-Null

# Row 34
$ws.Cells.Item(34,1).NumberFormat = "@"
$ws.Cells.Item(34,1).Value = "3.4.2020"
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34,2).Value = 0.39583333333333331
$ws.Cells.Item(34,3).Value = 0.43055555555555558
$ws.Cells.Item(34,4).Formula = "=C34-B34"
$ws.Cells.Item(34,5).Value = "Clk PLL"
$ws.Cells.Item(34,6).Value = "Create necessary files"
$ws.Cells.Item(34,7).Value = "IP Core"

# Row 35
$ws.Cells.Item(35,1).NumberFormat = "@"
$ws.Cells.Item(35,1).Value = "3.4.2020"
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35,2).Value = 0.43055555555555558
$ws.Cells.Item(35,3).Value = 0.45833333333333331
$ws.Cells.Item(35,4).Formula = "=C35-B35"
$ws.Cells.Item(35,5).Value = "Clk PLL"
$ws.Cells.Item(35,6).Value = "Testbench"

$ws.Range("G35").Select() | Out-Null

Write-Host "done"
